# TariffDecision1 / GENERAL CONSULTATION row is being removed from the
# "Feuil1" rule table. Select the whole row 6 (as Excel does when a user
# right-clicks a row header) and delete it, shifting rows 7-9 up to 6-8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Rows(6)
$selected = $row.Select()
$row.Delete()
